$wb = $excel.ActiveWorkbook

$wsReg   = $wb.Sheets.Item(1)   # "Registration"
$wsLogin = $wb.Sheets.Item(2)   # "login"

# Update the email address used on both sheets (shared string).
# Updating both occurrences lets the engine collapse the old/new
# shared-string entries so every reference stays in sync.
$wsReg.Range("D2").Value = "saivara34256@gmail.com"
$wsLogin.Range("A2").Value = "saivara34256@gmail.com"

# Registration sheet: no longer the selected tab, new active cell D2.
$wsReg.Activate()
$wsReg.Range("D2").Select()

# login sheet: becomes the selected tab, new active cell B7.
$wsLogin.Activate()
$wsLogin.Range("B7").Select()
